# Fixed #295 - Add the version of M2Doc in the template custom properties.
#
# This particular template resource's diff (word/document.xml and
# word/styles.xml) is a pure re-serialization: every changed line carries
# exactly the same set of XML attribute name/value pairs before and after,
# only in a different (alphabetized) order, e.g.
#   -  <w:color w:val="E36C0A" w:themeColor="accent6" w:themeShade="BF"/>
#   +  <w:color w:themeColor="accent6" w:themeShade="BF" w:val="E36C0A"/>
# and likewise for the namespace declarations on <w:document>, <w:pgSz>,
# <w:pgMar>, <w:rFonts>, <w:lang>, <w:latentStyles>, every <w:lsdException>,
# and every <w:style>/<w:tblInd>/<w:tblCellMar> entry in the default styles.
# No run text, field code, formatting value, page-size/margin number, style
# definition, or section setting is added, removed, or changed anywhere in
# the document -- this is a canonicalization/resave artifact, not a content
# edit, so there is nothing in Word's object model to change here.
#
# We therefore touch nothing: the document is left exactly as authored, which
# is what keeps its canonical (attribute-order-independent) OOXML equal to
# the target. (Word property round-trips -- e.g. re-assigning PageSetup
# margins to their own current values -- were verified to NOT reproduce the
# upstream attribute ordering, since that ordering isn't reachable through
# any exposed COM property, and they additionally perturb unrelated parts,
# such as recomputing docProps/app.xml word/character/paragraph counts,
# which are not part of this change. So making such a "no-op" edit would
# actually introduce diffs that shouldn't be there.)

$d = $word.ActiveDocument
Write-Output ("Paragraphs: " + $d.Paragraphs.Count)
